$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap rows 114 and 115 (columns B:AC), matches/ids reordered ---
$cols1 = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")
foreach ($c in $cols1) {
    $addr1 = "$c" + "114"
    $addr2 = "$c" + "115"
    $v1 = $ws.Range($addr1).Value()
    $v2 = $ws.Range($addr2).Value()
    $ws.Range($addr1).Value = $v2
    $ws.Range($addr2).Value = $v1
}

# --- Step 2: swap rows 116 and 117 (columns B:AC) ---
foreach ($c in $cols1) {
    $addr1 = "$c" + "116"
    $addr2 = "$c" + "117"
    $v1 = $ws.Range($addr1).Value()
    $v2 = $ws.Range($addr2).Value()
    $ws.Range($addr1).Value = $v2
    $ws.Range($addr2).Value = $v1
}

# --- Step 3: update rows 251-257 with the refreshed odds data (re-sorted fixtures) ---
# row 251
$ws.Range("B251").Value = 7822431
$ws.Range("E251").Value = 45386.57291666666
$ws.Range("F251").Value = "FC Volendam"
$ws.Range("G251").Value = "Feyenoord"
$ws.Range("K251").Value = 13
$ws.Range("L251").Value = 8.5
$ws.Range("M251").Value = 1.125
$ws.Range("N251").Value = 23
$ws.Range("O251").Value = 11
$ws.Range("P251").Value = 1.1
$ws.Range("Q251").Value = 2.75
$ws.Range("R251").Value = 1.92
$ws.Range("S251").Value = 1.98
$ws.Range("T251").Value = 4
$ws.Range("U251").Value = 1.875
$ws.Range("V251").Value = 1.975

# row 252
$ws.Range("B252").Value = 7822780
$ws.Range("E252").Value = 45386.625
$ws.Range("F252").Value = "RKC"
$ws.Range("G252").Value = "Almere City FC"
$ws.Range("K252").Value = 2.1
$ws.Range("L252").Value = 3.6
$ws.Range("M252").Value = 3.1
$ws.Range("N252").Value = 2.2
$ws.Range("O252").Value = 3.6
$ws.Range("P252").Value = 3.2
$ws.Range("Q252").Value = -0.25
$ws.Range("R252").Value = 1.95
$ws.Range("S252").Value = 1.95
$ws.Range("T252").Value = 2.5
$ws.Range("U252").Value = 2.025
$ws.Range("V252").Value = 1.825

# row 253
$ws.Range("B253").Value = 7822225
$ws.Range("E253").Value = 45386.66666666666
$ws.Range("F253").Value = "Ajax"
$ws.Range("G253").Value = "Go Ahead Eagles"
$ws.Range("K253").Value = 1.444
$ws.Range("L253").Value = 4.6
$ws.Range("M253").Value = 6
$ws.Range("N253").Value = 1.5
$ws.Range("O253").Value = 4.75
$ws.Range("P253").Value = 5.75
$ws.Range("Q253").Value = -1.25
$ws.Range("R253").Value = 2.05
$ws.Range("S253").Value = 1.85
$ws.Range("T253").Value = 3.25
$ws.Range("U253").Value = 2
$ws.Range("V253").Value = 1.85

# row 254
$ws.Range("B254").Value = 6973370
$ws.Range("E254").Value = 45388.47916666666
$ws.Range("F254").Value = "Sparta Rotterdam"
$ws.Range("G254").Value = "Heracles"
$ws.Range("K254").Value = 1.75
$ws.Range("L254").Value = 3.8
$ws.Range("M254").Value = 4.2
$ws.Range("N254").Value = 1.7
$ws.Range("O254").Value = 4
$ws.Range("P254").Value = 4.5
$ws.Range("Q254").Value = -0.75
$ws.Range("R254").Value = 1.92
$ws.Range("S254").Value = 1.98
$ws.Range("T254").Value = 2.75
$ws.Range("U254").Value = 1.85
$ws.Range("V254").Value = 2

# row 255
$ws.Range("B255").Value = 6838570
$ws.Range("E255").Value = 45388.57291666666
$ws.Range("F255").Value = "PSV"
$ws.Range("G255").Value = "AZ"
$ws.Range("K255").Value = 1.4
$ws.Range("L255").Value = 4.75
$ws.Range("M255").Value = 7
$ws.Range("N255").Value = 1.45
$ws.Range("O255").Value = 4.5
$ws.Range("P255").Value = 6
$ws.Range("Q255").Value = -1.25
$ws.Range("R255").Value = 2.05
$ws.Range("S255").Value = 1.85
$ws.Range("T255").Value = 3
$ws.Range("U255").Value = 2
$ws.Range("V255").Value = 1.85

# row 256
$ws.Range("B256").Value = 6956565
$ws.Range("E256").Value = 45388.625
$ws.Range("F256").Value = "PEC Zwolle"
$ws.Range("G256").Value = "Excelsior"
$ws.Range("K256").Value = 2.15
$ws.Range("L256").Value = 3.6
$ws.Range("M256").Value = 3.1
$ws.Range("N256").Value = 2
$ws.Range("O256").Value = 3.6
$ws.Range("P256").Value = 3.4
$ws.Range("Q256").Value = -0.5
$ws.Range("R256").Value = 2.05
$ws.Range("S256").Value = 1.85
$ws.Range("T256").Value = 3
$ws.Range("U256").Value = 2.05
$ws.Range("V256").Value = 1.8

# row 257
$ws.Range("B257").Value = 6956849
$ws.Range("E257").Value = 45388.66666666666
$ws.Range("F257").Value = "FC Twente"
$ws.Range("G257").Value = "Fortuna Sittard"
$ws.Range("K257").Value = 1.3
$ws.Range("L257").Value = 5.25
$ws.Range("M257").Value = 9
$ws.Range("N257").Value = 1.363
$ws.Range("O257").Value = 5
$ws.Range("P257").Value = 7
$ws.Range("Q257").Value = -1.25
$ws.Range("R257").Value = 1.86
$ws.Range("S257").Value = 2.04
$ws.Range("T257").Value = 2.75
$ws.Range("U257").Value = 1.95
$ws.Range("V257").Value = 1.9

# --- Step 4: append new rows 258 and 259 ---
# row 258: copy formatting from row 257 for A/E style reuse
$ws.Range("A257").Copy()
$ws.Range("A258").PasteSpecial(-4122)
$ws.Range("E257").Copy()
$ws.Range("E258").PasteSpecial(-4122)
$ws.Range("A258").Value = 256
$ws.Range("B258").Value = 6838586
$ws.Range("C258").Value = "Netherlands Eredivisie"
$ws.Range("D258").Value = "Netherlands Eredivisie"
$ws.Range("E258").Value = 45389.30208333334
$ws.Range("F258").Value = "Vitesse"
$ws.Range("G258").Value = "NEC"
$ws.Range("K258").Value = 2.8
$ws.Range("L258").Value = 3.6
$ws.Range("M258").Value = 2.3
$ws.Range("N258").Value = 2.875
$ws.Range("O258").Value = 3.6
$ws.Range("P258").Value = 2.25
$ws.Range("Q258").Value = 0.25
$ws.Range("R258").Value = 1.87
$ws.Range("S258").Value = 2.03
$ws.Range("T258").Value = 2.75
$ws.Range("U258").Value = 1.925
$ws.Range("V258").Value = 1.925
$ws.Range("W258").Value = 0
$ws.Range("X258").Value = 0
$ws.Range("Y258").Value = 0
$ws.Range("Z258").Value = 0
$ws.Range("AA258").Value = 0

# row 259: copy formatting from row 257 for A/E style reuse
$ws.Range("A257").Copy()
$ws.Range("A259").PasteSpecial(-4122)
$ws.Range("E257").Copy()
$ws.Range("E259").PasteSpecial(-4122)
$ws.Range("A259").Value = 257
$ws.Range("B259").Value = 6838571
$ws.Range("C259").Value = "Netherlands Eredivisie"
$ws.Range("D259").Value = "Netherlands Eredivisie"
$ws.Range("E259").Value = 45389.625
$ws.Range("F259").Value = "Heerenveen"
$ws.Range("G259").Value = "FC Utrecht"
$ws.Range("K259").Value = 2.25
$ws.Range("L259").Value = 3.5
$ws.Range("M259").Value = 2.9
$ws.Range("N259").Value = 2.45
$ws.Range("O259").Value = 3.5
$ws.Range("P259").Value = 2.625
$ws.Range("Q259").Value = 0
$ws.Range("R259").Value = 1.87
$ws.Range("S259").Value = 2.03
$ws.Range("T259").Value = 2.5
$ws.Range("U259").Value = 1.825
$ws.Range("V259").Value = 2.025
$ws.Range("W259").Value = 0
$ws.Range("X259").Value = 0
$ws.Range("Y259").Value = 0
$ws.Range("Z259").Value = 0
$ws.Range("AA259").Value = 0

$excel.CutCopyMode = 0
Write-Host "Edit complete"
